# Update pinout so pin names match schematic (MEPOS V0.5)
# - Tilde suffix moved from prefix ("~GPIOx") to suffix ("GPIOx~")
# - Legend entry symbol changed from "(~)GPIOx" to "GPIOx(~)"
# - Selection moved from M18 to L20

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pinout")

# Legend (key) entries, row 5
$ws.Range("L5").Value = "GPIOx(~)"
$ws.Range("M5").Value = "General Purpose IO (Timer/PWM)"

# Pin table entries - Odd/"B" column pins
$ws.Range("B19").Value = "GPIO2~"
$ws.Range("B20").Value = "GPIO4~"
$ws.Range("B26").Value = "GPIO16~"
$ws.Range("B27").Value = "GPIO18~"
$ws.Range("B28").Value = "GPIO20~"

# Pin table entries - Even/"E" column pins.
# These cells use a quote-prefixed style (because their text begins with
# "~"); writing .Value resets that formatting, so re-apply the original
# cell format afterwards by copying it from an unmodified cell that shares
# the same style (E21 keeps its original formatting all along).
$ws.Range("E19").Value = "GPIO1~"
$ws.Range("E20").Value = "GPIO3~"
$ws.Range("E26").Value = "GPIO15~"
$ws.Range("E27").Value = "GPIO17~"
$ws.Range("E28").Value = "GPIO19~"

$ws.Range("E21").Copy()
$ws.Range("E19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E21").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E21").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E21").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E21").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active cell selection to match the saved view state
$ws.Activate()
$ws.Range("L20").Select()
